$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '58.766.01'
$ws.Range("E2").Value = '  +0.40%  '

# Row 3
$ws.Range("D3").Value = '2.625.44'
$ws.Range("E3").Value = '  +3.41%  '

# Row 4
$ws.Range("E4").Value = '  +0.26%  '

# Row 5
$ws.Range("D5").Value = "'515.46"
$ws.Range("E5").Value = '  +1.42%  '

# Row 6
$ws.Range("D6").Value = "'143.28"
$ws.Range("E6").Value = '  -0.69%  '

# Row 7
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = '  +0.07%  '

# Row 8
$ws.Range("D8").Value = "'0.564"
$ws.Range("E8").Value = '  -0.33%  '

# Row 9
$ws.Range("D9").Value = '2.649.09'
$ws.Range("E9").Value = '  +4.12%  '

# Row 10
$ws.Range("E10").Value = '  +0.83%  '

# Row 11
$ws.Range("D11").Value = "'0.104"
$ws.Range("E11").Value = '  +2.34%  '

# Row 12
$ws.Range("D12").Value = "'0.334"
$ws.Range("E12").Value = '  +0.91%  '

# Row 13
$ws.Range("E13").Value = '  -1.77%  '

# Row 14
$ws.Range("D14").Value = '3.111.69'
$ws.Range("E14").Value = '  +4.28%  '

# Row 15
$ws.Range("D15").Value = '58.787.70'
$ws.Range("E15").Value = '  +0.44%  '

# Row 16
$ws.Range("D16").Value = "'20.75"
$ws.Range("E16").Value = '  -0.02%  '

# Row 17
$ws.Range("D17").Value = "'0.0000136"
$ws.Range("E17").Value = '  +0.76%  '

# Row 18
$ws.Range("D18").Value = '2.644.01'
$ws.Range("E18").Value = '  +3.95%  '

# Row 19
$ws.Range("D19").Value = "'346.61"
$ws.Range("E19").Value = '  +3.42%  '

# Row 20
$ws.Range("D20").Value = "'4.51"
$ws.Range("E20").Value = '  -0.66%  '

# Row 21
$ws.Range("D21").Value = "'10.27"
$ws.Range("E21").Value = '  +1.60%  '

# Row 22
$ws.Range("D22").Value = "'6.14"
$ws.Range("E22").Value = '  +2.84%  '

# Row 23
$ws.Range("D23").Value = "'0.999"
$ws.Range("E23").Value = '  +0.03%  '

# Row 24
$ws.Range("D24").Value = "'61.67"
$ws.Range("E24").Value = '  +2.02%  '

# Row 25
$ws.Range("D25").Value = "'0.417"
$ws.Range("E25").Value = '  +1.71%  '

# Row 26
$ws.Range("D26").Value = "'0.996"
$ws.Range("E26").Value = '  -0.83%  '

# Row 27
$ws.Range("D27").Value = "'0.161"
$ws.Range("E27").Value = '  +0.57%  '

# Row 28
$ws.Range("D28").Value = '0.0₃0795'
$ws.Range("E28").Value = '  +0.93%  '

# Row 29
$ws.Range("D29").Value = "'7.06"
$ws.Range("E29").Value = '  +1.11%  '

# Row 30
$ws.Range("D30").Value = "'1.00"
$ws.Range("E30").Value = '  +0.04%  '

# Row 31
$ws.Range("D31").Value = "'6.24"
$ws.Range("E31").Value = '  +6.55%  '

# Row 32
$ws.Range("D32").Value = "'18.86"
$ws.Range("E32").Value = '  +1.63%  '

# Row 33
$ws.Range("E33").Value = '  +2.28%  '

# Row 34
$ws.Range("D34").Value = "'149.34"
$ws.Range("E34").Value = '  -0.06%  '

# Row 35
$ws.Range("D35").Value = "'0.964"
$ws.Range("E35").Value = '  +5.09%  '

# Row 36
$ws.Range("D36").Value = "'3.98"
$ws.Range("E36").Value = '  +1.83%  '

# Row 37
$ws.Range("E37").Value = '  +1.37%  '

# Row 38
$ws.Range("D38").Value = "'36.61"
$ws.Range("E38").Value = '  +1.64%  '

# Row 39
$ws.Range("D39").Value = "'0.834"
$ws.Range("E39").Value = '  +0.95%  '

# Row 40
$ws.Range("D40").Value = "'3.69"
$ws.Range("E40").Value = '  +4.62%  '

# Row 41
$ws.Range("E41").Value = '  +0.13%  '

# Row 42
$ws.Range("B42").Value = 'FirstDigitalUSD'
$ws.Range("C42").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D42").Value = "'0.995"
$ws.Range("E42").Value = '  -0.33%  '

# Row 43
$ws.Range("B43").Value = 'Bittensor'
$ws.Range("C43").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D43").Value = "'275.65"
$ws.Range("E43").Value = '  -3.12%  '

# Row 44
$ws.Range("B44").Value = 'Mantle'
$ws.Range("C44").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D44").Value = "'0.608"
$ws.Range("E44").Value = '  +0.95%  '

# Row 45
$ws.Range("D45").Value = "'0.0982"
$ws.Range("E45").Value = '  -1.81%  '

# Row 46
$ws.Range("D46").Value = "'19.48"
$ws.Range("E46").Value = '  +4.33%  '

# Row 47
$ws.Range("D47").Value = "'0.0527"
$ws.Range("E47").Value = '  -1.34%  '

# Row 48
$ws.Range("D48").Value = "'10.29"
$ws.Range("E48").Value = '  -0.06%  '

# Row 49
$ws.Range("D49").Value = "'0.0228"
$ws.Range("E49").Value = '  +0.82%  '

# Row 50
$ws.Range("D50").Value = '1.974.49'
$ws.Range("E50").Value = '  +3.76%  '

# Row 51
$ws.Range("D51").Value = "'4.56"
$ws.Range("E51").Value = '  +0.36%  '
